$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2109375
$ws.Range("C2").Value = 0.54296875
$ws.Range("J2").Value = 0.015625
$ws.Range("O2").Value = 0.00390625
$ws.Range("P2").Value = 0.1328125
$ws.Range("S2").Value = 0.09375
$ws.Range("C3").Value = 0.02127659574468085
$ws.Range("J3").Value = 0.02127659574468085
$ws.Range("P3").Value = 0.7943262411347518
$ws.Range("S3").Value = 0.1631205673758865
$ws.Range("J4").Value = 0.02564102564102564
$ws.Range("P4").Value = 0.6923076923076923
$ws.Range("S4").Value = 0.282051282051282
$ws.Range("B6").Value = 0.03831417624521073
$ws.Range("D6").Value = 0.01532567049808429
$ws.Range("F6").Value = 0.0842911877394636
$ws.Range("J6").Value = 0.2222222222222222
$ws.Range("O6").Value = 0.01532567049808429
$ws.Range("Q6").Value = 0.1724137931034483
$ws.Range("R6").Value = 0.08045977011494253
$ws.Range("S6").Value = 0.3716475095785441
$ws.Range("B7").Value = 0.1116751269035533
$ws.Range("D7").Value = 0.02030456852791878
$ws.Range("F7").Value = 0.08121827411167512
$ws.Range("J7").Value = 0.08629441624365482
$ws.Range("O7").Value = 0.01522842639593909
$ws.Range("Q7").Value = 0.1370558375634518
$ws.Range("R7").Value = 0.1370558375634518
$ws.Range("S7").Value = 0.4111675126903553
$ws.Range("B8").Value = 0.08743169398907104
$ws.Range("D8").Value = 0.009107468123861567
$ws.Range("E8").Value = 0.00546448087431694
$ws.Range("F8").Value = 0.05828779599271403
$ws.Range("J8").Value = 0.1038251366120219
$ws.Range("O8").Value = 0.007285974499089253
$ws.Range("Q8").Value = 0.1730418943533698
$ws.Range("R8").Value = 0.1147540983606557
$ws.Range("S8").Value = 0.4408014571948998
$ws.Range("B9").Value = 0.04591836734693878
$ws.Range("D9").Value = 0.03061224489795918
$ws.Range("F9").Value = 0.06122448979591837
$ws.Range("J9").Value = 0.09693877551020408
$ws.Range("O9").Value = 0.02551020408163265
$ws.Range("Q9").Value = 0.1887755102040816
$ws.Range("R9").Value = 0.09693877551020408
$ws.Range("S9").Value = 0.4540816326530612
$ws.Range("B10").Value = 0.091699604743083
$ws.Range("D10").Value = 0.01739130434782609
$ws.Range("E10").Value = 0.002371541501976285
$ws.Range("F10").Value = 0.08537549407114625
$ws.Range("J10").Value = 0.09407114624505929
$ws.Range("O10").Value = 0.02055335968379447
$ws.Range("Q10").Value = 0.2142292490118577
$ws.Range("R10").Value = 0.09960474308300395
$ws.Range("S10").Value = 0.374703557312253
$ws.Range("G11").Value = 0.1277955271565495
$ws.Range("J11").Value = 0.1277955271565495
$ws.Range("K11").Value = 0.1980830670926517
$ws.Range("L11").Value = 0.5303514376996805
$ws.Range("S11").Value = 0.01597444089456869
$ws.Range("G12").Value = 0.7251461988304093
$ws.Range("J12").Value = 0.2280701754385965
$ws.Range("K12").Value = 0.01169590643274854
$ws.Range("L12").Value = 0.02923976608187134
$ws.Range("S12").Value = 0.005847953216374269
$ws.Range("G13").Value = 0.7857142857142857
$ws.Range("J13").Value = 0.2142857142857143
$ws.Range("F15").Value = 0.02347417840375587
$ws.Range("H15").Value = 0.2065727699530517
$ws.Range("I15").Value = 0.07511737089201878
$ws.Range("J15").Value = 0.352112676056338
$ws.Range("K15").Value = 0.06572769953051644
$ws.Range("M15").Value = 0.01408450704225352
$ws.Range("O15").Value = 0.02347417840375587
$ws.Range("S15").Value = 0.2394366197183098
$ws.Range("F16").Value = 0.02339181286549707
$ws.Range("H16").Value = 0.152046783625731
$ws.Range("I16").Value = 0.08771929824561403
$ws.Range("J16").Value = 0.4093567251461988
$ws.Range("K16").Value = 0.1169590643274854
$ws.Range("M16").Value = 0.01754385964912281
$ws.Range("N16").Value = 0.005847953216374269
$ws.Range("O16").Value = 0.05847953216374269
$ws.Range("S16").Value = 0.1286549707602339
$ws.Range("F17").Value = 0.01273885350318471
$ws.Range("H17").Value = 0.1804670912951168
$ws.Range("I17").Value = 0.08917197452229299
$ws.Range("J17").Value = 0.445859872611465
$ws.Range("K17").Value = 0.09129511677282377
$ws.Range("M17").Value = 0.01061571125265393
$ws.Range("N17").Value = 0.004246284501061571
$ws.Range("O17").Value = 0.04246284501061571
$ws.Range("S17").Value = 0.1231422505307856
$ws.Range("F18").Value = 0.02745098039215686
$ws.Range("H18").Value = 0.2313725490196079
$ws.Range("I18").Value = 0.05098039215686274
$ws.Range("J18").Value = 0.4352941176470588
$ws.Range("K18").Value = 0.08235294117647059
$ws.Range("M18").Value = 0.01568627450980392
$ws.Range("N18").Value = 0.00392156862745098
$ws.Range("O18").Value = 0.05490196078431372
$ws.Range("S18").Value = 0.09803921568627451
$ws.Range("F19").Value = 0.01561338289962825
$ws.Range("H19").Value = 0.2490706319702602
$ws.Range("I19").Value = 0.08252788104089219
$ws.Range("J19").Value = 0.3405204460966543
$ws.Range("K19").Value = 0.1048327137546468
$ws.Range("M19").Value = 0.020817843866171
$ws.Range("N19").Value = 0.0007434944237918215
$ws.Range("O19").Value = 0.06617100371747212
$ws.Range("S19").Value = 0.1197026022304833
